# PM03 Tidsregistrering for Rasmus.xlsx - add new time-registration entry
# for row 21 on sheet "Ark1": a new task "rettelse af UCD01 efter samtale
# med Ander" performed on 2020-05-14 (serial 43965) from 08:30 to 10:00,
# and move the active selection to D22 (previously A22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# New log entry in row 21
$ws.Range("A21").Value = "rettelse af UCD01 efter samtale med Ander"
$ws.Range("C21").Value = 43965
$ws.Range("D21").Value = 0.35416666666666669
$ws.Range("E21").Value = 0.41666666666666669

# Move selection from A22 to D22
$ws.Activate()
$ws.Range("D22").Select()
